$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$rowCount = $t.Rows.Count
$colCount = $t.Columns.Count

# Light-grey (D3D3D3) 0.5pt single borders added to every cell of the table.
$borderColor = 13882323   # BGR int for D3D3D3

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $t.Cell($r, $c)
        foreach ($side in -1, -2, -3, -4) {
            $border = $cell.Borders.Item($side)
            $border.LineStyle = 1
            $border.LineWidth = 2
            $border.Color = $borderColor
        }
    }
}

# Header row (row 1) shading: FFEFD5 -> FFDAB9.
$headerFill = 12180223   # BGR int for FFDAB9
for ($c = 1; $c -le $colCount; $c++) {
    $headerCell = $t.Cell(1, $c)
    $headerCell.Shading.BackgroundPatternColor = $headerFill
}

# Data rows (2..N), columns 2..N: paragraph alignment right -> center.
for ($r = 2; $r -le $rowCount; $r++) {
    for ($c = 2; $c -le $colCount; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.ParagraphFormat.Alignment = 1
    }
}
